$d = $word.ActiveDocument

# Locate the paragraph that currently reads "晴，终于到了周末我要出去浪".
# This is the last paragraph in the document before the edit; the commit
# inserts two new diary entries (for 2022-7-1) right after it and also
# drops the stray eastAsia font-hint that was sitting on that paragraph's
# paragraph-mark run properties.
$targetIndex = -1
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    if ($d.Paragraphs.Item($i).Range.Text -like "*出去浪*") {
        $targetIndex = $i
    }
}

if ($targetIndex -eq -1) {
    throw "Could not locate target paragraph containing 出去浪"
}

$targetPara = $d.Paragraphs.Item($targetIndex)

if ($targetIndex -gt 1) {
    $rangeStart = $d.Paragraphs.Item($targetIndex - 1).Range.End
} else {
    $rangeStart = $targetPara.Range.Start
}

$replaceRange = $d.Range($rangeStart, $d.Content.End)

$xml = @'
<?xml version="1.0" encoding="UTF-8" standalone="yes"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
  <pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">
    <pkg:xmlData>
      <w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
        <w:body>
          <w:p>
            <w:pPr>
              <w:rPr>
                <w:sz w:val="32"/>
                <w:szCs w:val="32"/>
              </w:rPr>
            </w:pPr>
            <w:r>
              <w:rPr>
                <w:rFonts w:hint="eastAsia"/>
                <w:sz w:val="32"/>
                <w:szCs w:val="32"/>
              </w:rPr>
              <w:t>晴，终于到了周末我要出去浪</w:t>
            </w:r>
          </w:p>
          <w:p>
            <w:pPr>
              <w:rPr>
                <w:sz w:val="32"/>
                <w:szCs w:val="32"/>
              </w:rPr>
            </w:pPr>
            <w:r>
              <w:rPr>
                <w:rFonts w:hint="eastAsia"/>
                <w:sz w:val="32"/>
                <w:szCs w:val="32"/>
              </w:rPr>
              <w:t>2</w:t>
            </w:r>
            <w:r>
              <w:rPr>
                <w:sz w:val="32"/>
                <w:szCs w:val="32"/>
              </w:rPr>
              <w:t>022-7-1</w:t>
            </w:r>
            <w:r>
              <w:rPr>
                <w:rFonts w:hint="eastAsia"/>
                <w:sz w:val="32"/>
                <w:szCs w:val="32"/>
              </w:rPr>
              <w:t>日星期五</w:t>
            </w:r>
          </w:p>
          <w:p>
            <w:pPr>
              <w:rPr>
                <w:rFonts w:hint="eastAsia"/>
                <w:sz w:val="32"/>
                <w:szCs w:val="32"/>
              </w:rPr>
            </w:pPr>
            <w:r>
              <w:rPr>
                <w:rFonts w:hint="eastAsia"/>
                <w:sz w:val="32"/>
                <w:szCs w:val="32"/>
              </w:rPr>
              <w:t>晴，今天天气不错奖励自己一个冰棒</w:t>
            </w:r>
            <w:r>
              <w:rPr>
                <w:rFonts w:hint="eastAsia"/>
                <w:sz w:val="32"/>
                <w:szCs w:val="32"/>
              </w:rPr>
              <w:t>，</w:t>
            </w:r>
            <w:r>
              <w:rPr>
                <w:rFonts w:hint="eastAsia"/>
                <w:sz w:val="32"/>
                <w:szCs w:val="32"/>
              </w:rPr>
              <w:t>心情真好</w:t>
            </w:r>
          </w:p>
        </w:body>
      </w:document>
    </pkg:xmlData>
  </pkg:part>
</pkg:package>
'@

$replaceRange.InsertXML($xml)
